$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (Total classes split across two dates) - E10/F10 were blank, now hold counts
$ws.Range("E10").Value = 18
$ws.Range("F10").Value = 19

# Attendance entries for the two new class dates (columns E and F) for each student row
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 3

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 3

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 3

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 3

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 3

$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 3

$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 3

$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 3

$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0

$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 0

$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0

$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0

$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 3

# Update the active selection to match the merged course-title cell range
$ws.Range("D8:H8").Select()
